$d = $word.ActiveDocument

# Locate the paragraph that currently reads:
#   "Template missing required/referenced parameter definition in parameter
#    section - fixed by adding missing parameters"
# and insert five new paragraphs immediately after it (i.e. right before the
# following empty "ListParagraph" paragraph), matching the target diff.

$findRange = $d.Content
$found = $findRange.Find.Execute(
    "Template missing required/referenced parameter definition in parameter section",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0
)

# Expand the found range to the full paragraph (includes its paragraph mark).
$findRange.Expand(4) | Out-Null
$target = $findRange

# Collapse to a point immediately before the anchor paragraph's own
# paragraph mark (i.e. the very end of its text), so the insertion lands
# between this paragraph and the following (untouched) empty paragraph.
$insertPoint = $d.Range($target.End - 1, $target.End - 1)

$w = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

$para1 = "<w:p $w><w:pPr><w:pStyle w:val='ListParagraph'/></w:pPr></w:p>"

$para2 = "<w:p $w><w:pPr><w:pStyle w:val='ListParagraph'/><w:numPr><w:ilvl w:val='0'/><w:numId w:val='2'/></w:numPr></w:pPr><w:r><w:t>17/08/2020, 10:37:37 - Template contains errors.: [/Outputs] 'null' values are not allowed in templates</w:t></w:r></w:p>"

$para3 = "<w:p $w>" +
    "<w:pPr><w:pStyle w:val='ListParagraph'/></w:pPr>" +
    "<w:r><w:t xml:space='preserve'>Fixed by outputting </w:t></w:r>" +
    "<w:proofErr w:type='spellStart'/><w:r><w:t>rLambdaFunction</w:t></w:r><w:proofErr w:type='spellEnd'/>" +
    "<w:r><w:t xml:space='preserve'> Alias and </w:t></w:r>" +
    "<w:proofErr w:type='spellStart'/><w:r><w:t>SSMParameter</w:t></w:r><w:proofErr w:type='spellEnd'/>" +
    "<w:r><w:t xml:space='preserve'> </w:t></w:r>" +
    "<w:proofErr w:type='spellStart'/><w:r><w:t>arn</w:t></w:r><w:proofErr w:type='spellEnd'/>" +
    "<w:r><w:t xml:space='preserve'> output</w:t></w:r>" +
    "</w:p>"

$para4 = "<w:p $w><w:pPr><w:pStyle w:val='ListParagraph'/></w:pPr></w:p>"
$para5 = "<w:p $w><w:pPr><w:pStyle w:val='ListParagraph'/></w:pPr></w:p>"

$xml = $para1 + $para2 + $para3 + $para4 + $para5

$insertPoint.InsertXML($xml)
